$d = $word.ActiveDocument

# ============================================================
# Helper: insert text at a collapsed range as a NEW run (not
# merged into the neighbouring run) by toggling Bold on then
# off again - the engine only coalesces adjacent runs with
# identical formatting during the insert itself, so a momentary
# difference forces a run boundary that survives even after the
# formatting is reverted.
# ============================================================
function Insert-AsNewRun {
    param($range, [string]$text)
    $range.InsertBefore($text)
    $ins = $d.Range($range.Start, $range.Start + $text.Length)
    $ins.Bold = 1
    $ins.Bold = 0
}

# ============================================================
# Hunk 1: "Руководитель ..." paragraph - replace placeholder
# blank with the filled-in supervisor details.
# ============================================================
$para1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Власов Дмитрий*") {
        $para1 = $p
        break
    }
}
$p1Start = $para1.Range.Start
$p1End = $para1.Range.End

# a) shrink the trailing 70-underscore run down to 15 underscores
$r = $d.Range($p1Start, $p1End)
$r.Find.Execute("______________________________________________________________________", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = "_______________"

# b) drop " электронного обучения"
$r = $d.Range($p1Start, $p1End)
$r.Find.Execute(" электронного обучения", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = ""

# c) "информационных технологий и" -> "ИТиЭО" (kept as its own run,
#    mirroring the proofErr-wrapped run the real diff introduces)
$r = $d.Range($p1Start, $p1End)
$r.Find.Execute("информационных технологий и", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = ""
Insert-AsNewRun $r "ИТиЭО"

# d) insert a new (non-underlined) 15-underscore run right after "Руководитель "
$r = $d.Range($p1Start, $p1End)
$r.Find.Execute("Руководитель ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
Insert-AsNewRun $r "_______________"

Write-Host "Hunk 1 done"
